# Fruta / hortaliza, semanal
# The weekly refresh re-shuffles the per-record fields (Fecha, Variedad,
# Calidad, Volumen, Precio minimo/maximo/promedio, Unidad de
# comercializacion, Origen, Precio $/Kg) across the existing data rows of
# the sheet. Column A/B/C/E/F/G/Q/R stay constant for every row, so the
# edit only ever touches D, H, I, J, K, L, M, N, O, P for rows 2..22.
#
# Snapshot the "before" values for those columns, then write them back out
# under the new row order described by $map (new row -> source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 22
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P")

# new row number -> source row number (values picked up from that source
# row get written into the new row)
$map = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 17
    6  = 15
    7  = 16
    8  = 4
    9  = 8
    10 = 10
    11 = 18
    12 = 19
    13 = 20
    14 = 21
    15 = 9
    16 = 12
    17 = 3
    18 = 22
    19 = 2
    20 = 13
    21 = 14
    22 = 11
}

# 1) Snapshot existing values per row/column before any writes happen.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each row's cells using the snapshotted values from its mapped
#    source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $map[$r]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
